# Update expiry dates (column B) to push subscriptions out by roughly
# 4 weeks, so only a single reminder email goes out per cycle.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44993
$ws.Range("B3").Value = 45000
$ws.Range("B4").Value = 45001
$ws.Range("B5").Value = 45000

# Move the active selection to B5, matching where the last edit was made.
$ws.Range("B5").Select()
